$d = $word.ActiveDocument

$old = "1)Telles, P. C. S. - Materiais para Equipamentos de Processo - Ed. Interciência, 4º Ed., 1989.2)Bresciani, F., E. - Seleção de Materiais Metálicos - Ed. da UNICAMP, 2º Ed.3)Freire, J. M. -Materiais de Construção Mecânica - Ed. Livros Técnicos e Científicos, Editora 1993.4)A. Remy/ M. Gay/ R. Gonthier - Materiais - Hemus Editora Limitada - 2ª Edição.5)Chiaverini, V.Tecnologia Mecânica - Materiais de Construção Mecânica - Vol. II - Ed. McGraw Hill do Brasil Ltda.6)Gentil, V. - Corrosão. - Ed. Guanabara Dois, 1982."

$new = "1)Telles, P. C. S. - Materiais para Equipamentos de Processo - Ed. Interciência, 4º Ed., 1989.^l2)Bresciani, F., E. - Seleção de Materiais Metálicos - Ed. da UNICAMP, 2º Ed.^l3)Freire, J. M. -Materiais de Construção Mecânica - Ed. Livros Técnicos e Científicos, Editora 1993.^l4)A. Remy/ M. Gay/ R. Gonthier - Materiais - Hemus Editora Limitada - 2ª Edição.^l5)Chiaverini, V.Tecnologia Mecânica - Materiais de Construção Mecânica - Vol. II - Ed. McGraw Hill do Brasil Ltda.^l6)Gentil, V. - Corrosão. - Ed. Guanabara Dois, 1982."

$found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)

Write-Host "Found and replaced: $found"
